# Added Dictionary for Every Test Case
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename Sheet2 -> ProductTestData and populate with product test data
# ---------------------------------------------------------------------------
$product = $wb.Worksheets.Item("Sheet2")
$product.Name = "ProductTestData"

$productData = @(
    @("TestCaseName",            "Execution Required", "keyword", "brand"),
    @("relevantProducts",        "Yes ",                "iphone", "APPLE"),
    @("sortProductByPriceHtoL",  "Yes ",                "iphone", $null),
    @("sortProductByPriceLtoH",  "Yes ",                "iphone", $null),
    @("chooseProductBrand",      "Yes ",                "mobile", "apple"),
    @("chooseProductRating",     "Yes ",                "iphone", $null)
)

for ($r = 0; $r -lt $productData.Count; $r++) {
    $row = $productData[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        if ($null -ne $row[$c]) {
            $product.Cells.Item($r + 1, $c + 1).Value = $row[$c]
        }
    }
}

$product.Range("A1:D1").Font.Bold = $true

$product.Columns.Item(1).ColumnWidth = 24.166666666666668
$product.Columns.Item(2).ColumnWidth = 17.893229166666668
$product.Columns.Item(3).ColumnWidth = 23.256510416666668

$product.Range("A1:D3").Select()

# ---------------------------------------------------------------------------
# 2. Add WishlistTestData after ProductTestData
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wishlist = $wb.Worksheets.Add($null, $lastSheet)
$wishlist.Name = "WishlistTestData"

$wishlistHeaders = @("TestCaseName", "Execution Required", "keyword", "brand", "mob", "pwd", "keyword")
for ($c = 0; $c -lt $wishlistHeaders.Count; $c++) {
    $wishlist.Cells.Item(1, $c + 1).Value = $wishlistHeaders[$c]
}

$wishlistData = @(
    @("addTowishlist",        "No",  "iphone", "APPLE", '"8708185463"', "pulkit9017", "iPhone"),
    @("removeFromwishlist",   "No",  "iphone", $null,   '"8708185463"', "pulkit9017", $null),
    @("getProductsInwishlist","Yes", "iphone", $null,   '"8708185463"', "pulkit9017", $null)
)

for ($r = 0; $r -lt $wishlistData.Count; $r++) {
    $row = $wishlistData[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        if ($null -ne $row[$c]) {
            $wishlist.Cells.Item($r + 2, $c + 1).Value = $row[$c]
        }
    }
}

$wishlist.Range("A1:G1").Font.Bold = $true

$wishlist.Columns.Item(1).ColumnWidth = 21.709635416666668
$wishlist.Columns.Item(2).ColumnWidth = 22.346354166666668
$wishlist.Columns.Item(5).ColumnWidth = 9.983072916666666
$wishlist.Columns.Item(6).ColumnWidth = 12.983072916666666

$wishlist.Activate()
$wishlist.Range("B3").Select()

# ---------------------------------------------------------------------------
# 3. Add the remaining (currently empty) test-data sheets
# ---------------------------------------------------------------------------
$sheetNames = @(
    "AddressTestData",
    "LoginTestData",
    "HomePageTestData",
    "LogoutTestData",
    "CartTestData",
    "ProductDetailsTestData",
    "FlipkartProfileInfoTestData"
)

foreach ($name in $sheetNames) {
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add($null, $last)
    $newSheet.Name = $name
}

$wb.Worksheets.Item("ProductDetailsTestData").Range("K17").Select()
$wb.Worksheets.Item("FlipkartProfileInfoTestData").Range("H14").Select()

# ---------------------------------------------------------------------------
# 4. Sheet1 selection update + re-activate WishlistTestData as the active tab
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("C2").Select()

$wishlist.Activate()
$wishlist.Range("B3").Select()
